$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 223, shifting existing rows 223:304 down to 224:305
$ws.Rows.Item(223).Insert()

# Populate the new row 223 with the new data record
$ws.Cells.Item(223, 1).Value = 9
$ws.Cells.Item(223, 2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(223, 3).Value = 'Metropolitana'
$ws.Cells.Item(223, 4).Value = 44809
$ws.Cells.Item(223, 5).Value = 13
$ws.Cells.Item(223, 6).Value = 100112001
$ws.Cells.Item(223, 7).Value = 'Berenjena'
$ws.Cells.Item(223, 8).Value = 'Sin especificar'
$ws.Cells.Item(223, 9).Value = 'Primera'
$ws.Cells.Item(223, 10).Value = 160
$ws.Cells.Item(223, 11).Value = 10000
$ws.Cells.Item(223, 12).Value = 12000
$ws.Cells.Item(223, 13).Value = 10875
$ws.Cells.Item(223, 14).Value = '$/caja 50 unidades'
$ws.Cells.Item(223, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(223, 16).Value = 218
$ws.Cells.Item(223, 17).Value = 50
$ws.Cells.Item(223, 18).Value = 'Hortaliza'
